# Swap the order of the two comma-separated names/emails in column G
# ("Recorded By") wherever the value ends with ", dnasr281@gmail.com"
# (i.e. "X, dnasr281@gmail.com" -> "dnasr281@gmail.com, X").
# Cells already starting with "dnasr281@gmail.com" or containing only a
# single value are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$suffix = ", dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.EndsWith($suffix)) {
        $other = $val.Substring(0, $val.Length - $suffix.Length)
        $cell.Value2 = "dnasr281@gmail.com, " + $other
    }
}
